$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.8268259657126666
$ws.Range("R2").Value = 7.441433691414001
$ws.Range("S2").Value = 0.0001592253407520804
$ws.Range("T2").Value = 0.0001592253407520804
$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("Q3").Value = 46.26664859026055
$ws.Range("R3").Value = 416.399837312345
$ws.Range("S3").Value = 0.008909762383782065
$ws.Range("T3").Value = 0.008909762383782067
$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 113.6783013069767
$ws.Range("R4").Value = 1023.10471176279
$ws.Range("S4").Value = 0.02189150681318974
$ws.Range("T4").Value = 0.02189150681318974
$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 15.98613895620133
$ws.Range("R5").Value = 143.875250605812
$ws.Range("S5").Value = 0.003078517763308638
$ws.Range("T5").Value = 0.003078517763308639
$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("Q6").Value = 894.5353727058344
$ws.Range("R6").Value = 8050.818354352509
$ws.Range("S6").Value = 0.172264424970143
$ws.Range("T6").Value = 0.172264424970143
$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 2197.895562498313
$ws.Range("R7").Value = 19781.06006248482
$ws.Range("S7").Value = 0.4232579580088991
$ws.Range("T7").Value = 0.4232579580088991
$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("M8").Value = 0.303146
$ws.Range("N8").Value = 0.909438
$ws.Range("O8").Value = 0.005142855213700541
$ws.Range("P8").Value = 0.005142855213700542
$ws.Range("Q8").Value = 9.892873536358001
$ws.Range("R8").Value = 89.035861827222
$ws.Range("S8").Value = 0.001905112109639822
$ws.Range("T8").Value = 0.001905112109639822
$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("O9").Value = 0.2877784259203595
$ws.Range("P9").Value = 0.2877784259203595
$ws.Range("Q9").Value = 553.5749026217984
$ws.Range("R9").Value = 4982.174123596185
$ws.Range("S9").Value = 0.1066042385664343
$ws.Range("T9").Value = 0.1066042385664343
$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 41.67881
$ws.Range("N10").Value = 125.03643
$ws.Range("O10").Value = 0.7070787188659401
$ws.Range("P10").Value = 0.7070787188659401
$ws.Range("Q10").Value = 1360.14724415263
$ws.Range("R10").Value = 12241.32519737367
$ws.Range("S10").Value = 0.2619292540438511
$ws.Range("T10").Value = 0.2619292540438511
